# edit.ps1 - applies the "Add files via upload" revision to Work Report.docx
#
# Summary of the target change (from the canonical-OOXML diff):
#   1. In the "Integration of ... flow for the user ..." bullet, the words
#      ", dashboard" are removed (the sentence used to read
#      "...add store, dashboard flow for the user..." and becomes
#      "...add store flow for the user...").
#   2. The whole bullet describing the now-removed
#      "webapp/api/dashboard/store-details/" GET route is deleted entirely
#      (it sat between the "webapp/api/signup/" POST bullet and the
#      "webapp/api/dashboard/add-store/" POST bullet).
#   3. styles.xml marks the built-in "Default Paragraph Font" character
#      style as semiHidden. This particular flag is not reachable through
#      the Style object's exposed COM surface in this host (Style.Hidden -
#      the VBA property that normally maps to <w:semiHidden/> - has no
#      working setter here, and the only wired-up visibility property,
#      Style.Visibility, writes <w:hidden/> instead of <w:semiHidden/>,
#      which would add an element the target document does not have), so
#      it is intentionally left untouched rather than writing incorrect
#      markup.

$d = $word.ActiveDocument

# --- Change 1: drop ", dashboard" from the integration-flow bullet -------
$found = $d.Content.Find.Execute(
    ", dashboard",  # FindText
    $true,          # MatchCase
    $false,         # MatchWholeWord
    $false,         # MatchWildcards
    $false,         # MatchSoundsLike
    $false,         # MatchAllWordForms
    $true,          # Forward
    1,              # Wrap (wdFindContinue)
    $false,         # Format
    "",             # ReplaceWith
    2               # Replace (wdReplaceAll)
)
if (-not $found) {
    Write-Output "WARNING: ', dashboard' text was not found (already edited?)"
}

# --- Change 2: remove the "dashboard/store-details" GET route bullet -----
$deleted = $false
for ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {
    $p = $d.Paragraphs($i)
    $t = $p.Range.Text
    if ($t -like "*store-details*") {
        $p.Range.Delete()
        $deleted = $true
        break
    }
}
if (-not $deleted) {
    Write-Output "WARNING: 'store-details' bullet paragraph was not found (already removed?)"
}

Write-Output "edit.ps1 complete"
